$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D28").Value = "Docker에서 tmux 사용법"
$ws.Range("E28").Value = "https://ropiens.tistory.com/162"

$ws.Range("D37").Value = "[Paper Review] Oversmoothing & Disassortative Graphs"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1826&mod=document&pageid=1"
